$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''69.294.55'
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  +2.23%  '
$ws.Range("D3").Value = '''3.834.10'
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").Value = '''630.09'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +5.16%  '
$ws.Range("D6").Value = '''166.55'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +0.72%  '
$ws.Range("D7").Value = '''3.834.27'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  +1.03%  '
$ws.Range("D8").Value = '''0.999'
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("D9").Value = '''0.522'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +0.87%  '
$ws.Range("E10").Value = '  +1.95%  '
$ws.Range("D11").Value = '''0.456'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  +0.76%  '
$ws.Range("D12").Value = '''6.60'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  +2.01%  '
$ws.Range("D13").Value = '''0.0000252'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  +1.02%  '
$ws.Range("D14").Value = '''36.18'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("D15").Value = '''4.469.96'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +0.82%  '
$ws.Range("D16").Value = '''3.839.59'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("D17").Value = '''69.257.45'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +2.14%  '
$ws.Range("D18").Value = '''18.17'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("D19").Value = '''7.17'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +1.45%  '
$ws.Range("D21").Value = '''467.78'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +0.99%  '
$ws.Range("D22").Value = '''9.73'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -0.88%  '
$ws.Range("D23").Value = '''0.712'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +1.61%  '
$ws.Range("D24").Value = '''0.0000155'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  +5.25%  '
$ws.Range("D25").Value = '''83.92'
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  +1.31%  '
$ws.Range("D26").Value = '''12.06'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("D27").Value = '''2.17'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +2.97%  '
$ws.Range("D28").Value = '''10.10'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  +0.63%  '
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").Value = '''3.975.56'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  +0.76%  '
$ws.Range("E31").Value = '  -1.20%  '
$ws.Range("E32").Value = '  +1.35%  '
$ws.Range("D33").Value = '''7.33'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -1.79%  '
$ws.Range("D34").Value = '''29.34'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +0.34%  '
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '''9.12'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  +0.95%  '
$ws.Range("B36").Value = 'Binance-PegBSC-USD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D36").Value = '''1.00'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("E37").Value = '  +2.36%  '
$ws.Range("E38").Value = '  +7.60%  '
$ws.Range("D39").Value = '''3.44'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +6.25%  '
$ws.Range("D40").Value = '''5.95'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  +3.16%  '
$ws.Range("D41").Value = '''0.982'
$ws.Range("D41").Style = 'Normal'
$ws.Range("D42").Value = '''0.999'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D44").Value = '''1.46'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +5.62%  '
$ws.Range("D45").Value = '''0.301'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +0.56%  '
$ws.Range("D46").Value = '''155.11'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +2.90%  '
$ws.Range("D47").Value = '''46.97'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -1.24%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = '''8.49'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +1.60%  '
$ws.Range("B49").Value = 'Arweave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D49").Value = '''42.65'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -4.72%  '
$ws.Range("D50").Value = '''1.90'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +2.40%  '
$ws.Range("D51").Value = '''0.000277'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +11.78%  '
